$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row for "None" style above the existing "Bokeh" row,
# then strip the inherited header formatting from the new row.
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).ClearFormats()

# Fill the new row's values (set B2 before A2 so the shared-string table
# ends up with " " allocated before "None ").
$ws.Range("B2").Value = " "
$ws.Range("A2").Value = "None "
$ws.Range("C2").Value = "OpenPose"

# Widen column B so the (now longer) description text fits.
$ws.Columns.Item(2).ColumnWidth = 49.67

# Update the active selection to C2.
$ws.Range("C2").Select()
